$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Majid Al Qassimi"
$summary.Range("B4").Value = 1627.94
$summary.Range("B6").Value = 2827
$summary.Range("B7").Value = 20207
$summary.Range("B8").Value = -17380
$summary.Range("B9").Value = 0.14

# --- Assets sheet ---
$assets = $wb.Worksheets.Item("Assets")
$assets.Range("C2").Value = 2827
$assets.Range("C3").Value = 2827

# --- Liabilities sheet ---
$liab = $wb.Worksheets.Item("Liabilities")
# Remove the "Personal Loans" row (row 2); "Credit Cards" row and "TOTAL LIABILITIES"
# row shift up to rows 2 and 3 respectively, keeping their original formatting.
$liab.Rows.Item(2).Delete()

# Update the (now shifted) Credit Cards row and TOTAL LIABILITIES row with new values.
$liab.Range("C2").Value = 20207
$liab.Range("D2").Value = 1010
$liab.Range("E2").Value = 1
$liab.Range("C3").Value = 20207
